$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) columns, and the three
# re-ordered rows (47-49) B/C/D/E cells, per the scraped refresh.

$ws.Range("D2").Value = "42.745.72"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.361.15"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'317.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.25%  "
$ws.Range("D6").Value = "'109.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("E7").Value = "  -2.24%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.76%  "
$ws.Range("D10").Value = "'42.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "'0.0926"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "'8.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  -4.66%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "'16.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.57%  "
$ws.Range("D16").Value = "2.718.61"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "2.254.34"
$ws.Range("E17").Value = "  -5.77%  "
$ws.Range("D18").Value = "42.725.30"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'7.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "'76.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "'3.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").Value = "'257.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.67%  "
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "'9.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'11.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").Value = "'22.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").Value = "'37.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "'172.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "'6.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").Value = "'2.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.73%  "
$ws.Range("E35").Value = "  +16.59%  "
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").Value = "'4.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").Value = "'3.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").Value = "'0.241"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("D42").Value = "'1.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.91%  "
$ws.Range("D43").Value = "'71.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").Value = "'112.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.77%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'5.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'9.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'86.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.21%  "
$ws.Range("D50").Value = "'77.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.90%  "
$ws.Range("E51").Value = "  -2.05%  "
